$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "magapoke_2025-11-19"

$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"

$hdr = $ws.Range("A1:B1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

$titles = @(
  'ブルーロック',
  'みいちゃんと山田さん',
  '信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐＆『ざまぁ！』します！',
  'WIND BREAKER',
  '東京卍リベンジャーズ',
  '宇宙兄弟',
  'ギルティサークル',
  'ドラハチ',
  '薫る花は凛と咲く',
  '島耕作',
  'ガチアクタ',
  '転生したら第七王子だったので、気ままに魔術を極めます',
  'イレギュラーズ',
  '十字架のろくにん',
  '君が僕らを悪魔と呼んだ頃',
  '黄昏町プリズナーズ',
  '黒猫と魔女の教室',
  '魔女と傭兵',
  'となりの黒川さん',
  'ハードワーカー中田',
  '魔術ギルド総帥～生まれ変わって今更やり直す2度目の学院生活～',
  '味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す',
  '転生貴族、鑑定スキルで成り上がる～弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた～',
  'K-9~警視庁公安部公安第9課異能対策係~',
  '異世界ウォーキング',
  '蒼く染めろ',
  'ひゃくえむ。',
  'アルキメデスの大戦',
  '幼馴染とはラブコメにならない',
  '追放された転生王子、『自動製作《オートクラフト》』スキルで領地を爆速で開拓し最強の村を作ってしまう〜最強クラフトスキルで始める、楽々領地開拓スローライフ〜',
  '愛妻の裏アカ',
  'FAIRY TAIL 100 YEARS QUEST',
  'グラぱらっ！',
  'アオバノバスケ',
  'いじめるヤバイ奴',
  '食糧人類-Starving Anonymous-',
  'せいぶつ部の田辺くん',
  '時々ボソッとロシア語でデレる隣のアーリャさん',
  'さわらないで小手指くん',
  'ともだちづくり',
  'ハナバス　苔石花江のバスケ論',
  'ストーカー行為がバレて人生終了男',
  'ナキナギ',
  '南海トラフ巨大地震',
  'しかのこのこのここしたんたん',
  'デッドアカウント',
  '辺境の薬師、都でSランク冒険者となる～英雄村の少年がチート薬で無自覚無双〜',
  'ジュミドロ',
  '阿武ノーマル',
  '念願の悪役令嬢（ラスボス）の身体を手に入れたぞ！',
  'Aランクパーティを離脱した俺は、元教え子たちと迷宮深部を目指す。',
  'この世界がいずれ滅ぶことを、俺だけが知っている～モンスターが現れた世界で、死に戻りレベルアップ～',
  '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～',
  '降り積もれ孤独な死よ',
  '屋根の下のアルテミス',
  '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～',
  '劣等人の魔剣使い　スキルボードを駆使して最強に至る',
  '春くらり',
  '四刀流の最強配信者～やり込んだVRゲームの設定が現実世界に反映されたので、廃止予定だった戦闘職で無双します～',
  'ダメスキル【自動機能】が覚醒しました～あれ、ギルドのスカウトの皆さん、俺を「いらない」って言ってませんでした？～',
  'おやすみ ふみさん',
  'なれの果ての僕ら',
  'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜',
  'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！',
  '不遇職【鑑定士】が実は最強だった～奈落で鍛えた最強の【神眼】で無双する～',
  'イジらないで、長瀞さん',
  '触手魔術師の成り上がり',
  'デスティニーラバーズ',
  '母という呪縛 娘という牢獄',
  'GALAXIAS',
  '東京卍リベンジャーズ～場地圭介からの手紙～',
  '恋ニ非ズ',
  'Destiny Unchain Online 〜吸血鬼少女となって、やがて『赤の魔王』と呼ばれるようになりました〜',
  'MYS',
  '田んぼで拾った女騎士、田舎で俺の嫁だと思われている',
  '我間乱 ―修羅―',
  'シャングリラ・フロンティア～クソゲーハンター、神ゲーに挑まんとす～',
  '復讐の教科書',
  'DAYS外伝',
  '現代転移の第二王子',
  'ヒロインは絶望しました。',
  '勇者と呼ばれた後に　―そして無双男は家族を創る―',
  '金田一少年の事件簿外伝 犯人たちの事件簿',
  'インフェクション',
  'お願い、脱がシて。',
  '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～',
  '追放されなかった男　～二度目の人生は土下座から始まりました～',
  '冰剣の魔術師が世界を統べる〜世界最強の魔術師である少年は、魔術学院に入学する〜',
  '英雄と魔女の転生ラブコメ',
  '最弱な僕は＜壁抜けバグ＞で成り上がる～壁をすり抜けたら、初回クリア報酬を無限回収できました！～',
  'お嬢様の僕',
  '可愛いだけじゃない式守さん',
  'それがメイドのカンナです',
  'ウイニング パス',
  '微妙に優しいいじめっ子',
  '絶望集落',
  '卒業アルバムの彼女たち',
  'はじめの一歩',
  '日本語が話せないロシア人美少女転入生が頼れるのは、多言語マスターの俺1人',
  '五輪の女神さま 〜なでしこ寮のメダルごはん〜'
)

for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $titles[$i]
}
